$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Addressing mode(Optional)" column (O) with a dropdown-style
# Static/Dynamic value set for rows 2 and 3 (寻址方式 switcher -> select).
# Shared-string insertion order matters for matching the target file, so
# write the header first, then row 3 (Static) before row 2 (Dynamic).
$ws.Range("O1").Value = "Addressing mode(Optional)"
$ws.Range("O3").Value = "Static"
$ws.Range("O2").Value = "Dynamic"

# Match the new column's width as closely as this runtime's pixel-quantized
# ColumnWidth model allows (target raw width is 27.125 char units).
$ws.Columns.Item(15).ColumnWidth = 26.4

# Restore the view: drop the custom top-left scroll position and move the
# active selection to G17 (matches the saved sheetView in the target file).
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("G17").Select()

# The "Speed <ZWSP><ZWSP>limit M/s(Optional)" header (N1) gets its two
# zero-width-space characters promoted to their own run in MS Gothic (the
# font Excel substitutes in for glyphs the primary font can't render) when
# Excel resaves this string. Recreate that run split.
$n1 = $ws.Range("N1")
$n1.Characters(7, 2).Font.Name = "MS Gothic"
